# Update "想去人数" (F column) and "最低票价" (G column) values for a handful
# of rows on both the "展览" and "全部类型" worksheets (they hold duplicate data).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# row number -> hashtable of column letter -> new value
$updates = @{
    3  = @{ F = 5008 }
    4  = @{ F = 5 }
    5  = @{ F = 7272 }
    12 = @{ F = 4252 }
    16 = @{ F = 2852 }
    19 = @{ F = 197 }
    20 = @{ F = 453 }
    21 = @{ F = 415 }
    22 = @{ F = 440 }
    23 = @{ F = 274 }
    24 = @{ F = 80 }
    28 = @{ F = 1333; G = 55 }
    36 = @{ F = 2668 }
    37 = @{ F = 686 }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $cols = $updates[$row]
        foreach ($col in $cols.Keys) {
            $ws.Range("$col$row").Value = $cols[$col]
        }
    }
}
